$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 120 (old rows 120-127 shift down to 122-129).
$ws.Range("A120:A121").EntireRow.Insert()

# --- New row 120 ---
$ws.Range("A120").Value2 = 8
$ws.Range("B120").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C120").Value2 = "Coquimbo"
$ws.Range("D120").Value2 = 44516
$ws.Range("D120").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E120").Value2 = 4
$ws.Range("F120").Value2 = 100112021
$ws.Range("G120").Value2 = "Ají"
$ws.Range("H120").Value2 = "Inferno"
$ws.Range("I120").Value2 = "Primera"
$ws.Range("J120").Value2 = 440
$ws.Range("K120").Value2 = 19000
$ws.Range("L120").Value2 = 20000
$ws.Range("M120").Value2 = 19500
$ws.Range("N120").Value2 = "$/caja 12 kilos"
$ws.Range("O120").Value2 = "Región de Arica y Parinacota"
$ws.Range("P120").Value2 = 1625
$ws.Range("Q120").Value2 = 12
$ws.Range("R120").Value2 = "Hortaliza"

# --- New row 121 ---
$ws.Range("A121").Value2 = 8
$ws.Range("B121").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C121").Value2 = "Coquimbo"
$ws.Range("D121").Value2 = 44516
$ws.Range("D121").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E121").Value2 = 4
$ws.Range("F121").Value2 = 100112021
$ws.Range("G121").Value2 = "Ají"
$ws.Range("H121").Value2 = "Inferno"
$ws.Range("I121").Value2 = "Segunda"
$ws.Range("J121").Value2 = 320
$ws.Range("K121").Value2 = 14000
$ws.Range("L121").Value2 = 15000
$ws.Range("M121").Value2 = 14500
$ws.Range("N121").Value2 = "$/caja 12 kilos"
$ws.Range("O121").Value2 = "Región de Arica y Parinacota"
$ws.Range("P121").Value2 = 1208
$ws.Range("Q121").Value2 = 12
$ws.Range("R121").Value2 = "Hortaliza"

$ws.Range("A1").Select()
